# Fruta / hortaliza, semanal
#
# Insert a new weekly price record for "Ciruela" (plum) variety "Automn Pride"
# at row 179 of Sheet1. Excel's row insert shifts every existing row from 179
# downward by one (old 179 -> 180, ..., old 201 -> 202), which is exactly the
# behaviour captured by the diff, and also grows the sheet's used range from
# A1:T201 to A1:T202 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 179..201 down to 180..202, leaving a blank row 179 to fill in.
$ws.Rows.Item(179).Insert()

$row = 179
$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value  = "Maule"
$ws.Cells.Item($row, 4).Value  = 44995
$ws.Cells.Item($row, 5).Value  = 7
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100103
$ws.Cells.Item($row, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value  = 100103002
$ws.Cells.Item($row, 10).Value = "Ciruela"
$ws.Cells.Item($row, 11).Value = "Automn Pride"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 820
$ws.Cells.Item($row, 14).Value = 9000
$ws.Cells.Item($row, 15).Value = 10000
$ws.Cells.Item($row, 16).Value = 9512
$ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row, 19).Value = 528
$ws.Cells.Item($row, 20).Value = 18
